$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Note: "17÷5=3, 2" is both a source value (row 5, col1) and a target value
# (row 2, col2 becomes "17÷5=3, 2"). We must replace the *existing*
# "17÷5=3, 2" occurrence BEFORE introducing the new one, otherwise the
# later replace-all would also catch the freshly written text.
Replace-Text "17÷5=3, 2" "78÷2=39, 0"

Replace-Text "75÷3=25, 0" "13÷9=1, 4"
Replace-Text "13÷3=4, 1" "44÷7=6, 2"
Replace-Text "40÷9=4, 4" "85÷5=17, 0"
Replace-Text "92÷4=23, 0" "42÷6=7, 0"
Replace-Text "11÷8=1, 3" "75÷2=37, 1"

Replace-Text "75÷7=10, 5" "19÷6=3, 1"
Replace-Text "21÷7=3, 0" "17÷5=3, 2"
Replace-Text "42÷5=8, 2" "76÷6=12, 4"
Replace-Text "45÷3=15, 0" "95÷8=11, 7"
Replace-Text "38÷5=7, 3" "81÷2=40, 1"

Replace-Text "89÷2=44, 1" "54÷7=7, 5"
Replace-Text "64÷2=32, 0" "70÷9=7, 7"
Replace-Text "94÷8=11, 6" "68÷4=17, 0"
Replace-Text "97÷7=13, 6" "42÷2=21, 0"
Replace-Text "50÷5=10, 0" "23÷5=4, 3"

Replace-Text "77÷3=25, 2" "74÷9=8, 2"
Replace-Text "59÷2=29, 1" "84÷7=12, 0"
Replace-Text "15÷2=7, 1" "14÷3=4, 2"
Replace-Text "37÷4=9, 1" "22÷6=3, 4"
Replace-Text "18÷5=3, 3" "64÷8=8, 0"

Replace-Text "14÷8=1, 6" "35÷9=3, 8"
Replace-Text "49÷6=8, 1" "80÷4=20, 0"
Replace-Text "26÷5=5, 1" "42÷9=4, 6"
Replace-Text "43÷4=10, 3" "14÷7=2, 0"
